$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.04"
$ws.Range("E2").Value = "'1.05%"
$ws.Range("D3").Value = "'30.45"
$ws.Range("E3").Value = "'11.50%"
$ws.Range("D4").Value = "'5.151"
$ws.Range("E4").Value = "'-0.12%"
$ws.Range("D5").Value = "'0.05730"
$ws.Range("E5").Value = "'1.45%"
$ws.Range("D6").Value = "'6.603"
$ws.Range("E6").Value = "'2.05%"
$ws.Range("D7").Value = "'3.047"
$ws.Range("E7").Value = "'1.44%"
$ws.Range("D8").Value = "'0.8598"
$ws.Range("E8").Value = "'4.91%"
$ws.Range("D9").Value = "'0.8725"
$ws.Range("E9").Value = "'4.39%"
$ws.Range("D10").Value = "'0.1364"
$ws.Range("E10").Value = "'2.42%"
$ws.Range("D11").Value = "'0.07098"
$ws.Range("E11").Value = "'2.47%"
$ws.Range("D12").Value = "'0.02866"
$ws.Range("E12").Value = "'-2.21%"
$ws.Range("D13").Value = "'0.09394"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.001517"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.04143"
$ws.Range("E15").Value = "'-3.13%"
$ws.Range("D16").Value = "'0.0005997"
$ws.Range("E16").Value = "'0.15%"
$ws.Range("D17").Value = "'0.005973"
$ws.Range("E17").Value = "'-2.49%"
$ws.Range("D18").Value = "'3.491"
$ws.Range("E18").Value = "'-0.48%"
$ws.Range("D19").Value = "'2.181"
$ws.Range("E19").Value = "'-5.64%"
$ws.Range("D20").Value = "'0.3196"
$ws.Range("E20").Value = "'2.66%"
$ws.Range("D21").Value = "'0.03241"
$ws.Range("E21").Value = "'4.23%"
$ws.Range("D22").Value = "'0.1308"
$ws.Range("E22").Value = "'1.22%"
$ws.Range("D23").Value = "'3.140"
$ws.Range("E23").Value = "'47.71%"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.005199"
$ws.Range("E25").Value = "'16.30%"
$ws.Range("B26").Value = "BitKan"
$ws.Range("C26").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D26").Value = "'0.001216"
$ws.Range("E26").Value = "'-0.75%"
$ws.Range("D27").Value = "'0.0001210"
$ws.Range("E27").Value = "'23.48%"
$ws.Range("D28").Value = "'0.0001374"
$ws.Range("E28").Value = "'89.39%"
$ws.Range("D40").Value = "'0.03778"
$ws.Range("E40").Value = "'3.60%"
$ws.Range("D41").Value = "'0.005779"
$ws.Range("E41").Value = "'-4.45%"
$ws.Range("E42").Value = "'1.88%"
$ws.Range("D43").Value = "'0.002600"
$ws.Range("E43").Value = "'13.06%"
$ws.Range("D44").Value = "'0.009782"
$ws.Range("E44").Value = "'19.24%"
$ws.Range("D45").Value = "'0.00005088"
$ws.Range("E45").Value = "'-5.06%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("D47").Value = "'0.07997"
$ws.Range("E47").Value = "'-20.80%"
$ws.Range("D48").Value = "'0.002771"
$ws.Range("E48").Value = "'3.98%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.02%"
